$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record ("Chino" / "Primera") was inserted before the
# existing row 179, pushing every subsequent record down by one row
# (old row 179 -> new row 180, ..., old row 272 -> new row 273).
$ws.Rows.Item(179).Insert()

$ws.Cells.Item(179, 1).Value = 8
$ws.Cells.Item(179, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(179, 3).Value = "Coquimbo"
$ws.Cells.Item(179, 4).Value = 44704
$ws.Cells.Item(179, 5).Value = 4
$ws.Cells.Item(179, 6).Value = 100112003
$ws.Cells.Item(179, 7).Value = "Ajo"
$ws.Cells.Item(179, 8).Value = "Chino"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 440
$ws.Cells.Item(179, 11).Value = 18500
$ws.Cells.Item(179, 12).Value = 19000
$ws.Cells.Item(179, 13).Value = 18750
$ws.Cells.Item(179, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(179, 15).Value = "China"
$ws.Cells.Item(179, 16).Value = 1875
$ws.Cells.Item(179, 17).Value = 10
$ws.Cells.Item(179, 18).Value = "Hortaliza"
